$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "B2" = 0.2293906810035842
    "C2" = 0.4982078853046595
    "J2" = 0.007168458781362007
    "P2" = 0.1863799283154122
    "S2" = 0.07885304659498207
    "B3" = 0.01398601398601399
    "C3" = 0.01398601398601399
    "J3" = 0.02097902097902098
    "P3" = 0.7272727272727273
    "S3" = 0.2237762237762238
    "J4" = 0.119047619047619
    "P4" = 0.6666666666666666
    "S4" = 0.2142857142857143
    "B6" = 0.0391304347826087
    "D6" = 0.004347826086956522
    "F6" = 0.05217391304347826
    "J6" = 0.2782608695652174
    "O6" = 0.04347826086956522
    "Q6" = 0.1521739130434783
    "R6" = 0.06086956521739131
    "S6" = 0.3695652173913043
    "B7" = 0.1126126126126126
    "D7" = 0.03153153153153153
    "F7" = 0.04054054054054054
    "J7" = 0.1981981981981982
    "O7" = 0.02702702702702703
    "Q7" = 0.1351351351351351
    "R7" = 0.07207207207207207
    "S7" = 0.3828828828828829
    "B8" = 0.0891566265060241
    "D8" = 0.01686746987951807
    "F8" = 0.06265060240963856
    "J8" = 0.09397590361445783
    "O8" = 0.02891566265060241
    "Q8" = 0.1493975903614458
    "R8" = 0.0963855421686747
    "S8" = 0.4626506024096386
    "B9" = 0.04864864864864865
    "D9" = 0.01621621621621622
    "F9" = 0.03783783783783784
    "J9" = 0.145945945945946
    "O9" = 0.02702702702702703
    "Q9" = 0.1945945945945946
    "R9" = 0.1351351351351351
    "S9" = 0.3945945945945946
    "B10" = 0.106280193236715
    "D10" = 0.02012882447665056
    "F10" = 0.06763285024154589
    "J10" = 0.1239935587761675
    "O10" = 0.02093397745571659
    "Q10" = 0.1948470209339775
    "R10" = 0.06763285024154589
    "S10" = 0.3985507246376812
    "G11" = 0.1820809248554913
    "J11" = 0.09826589595375723
    "K11" = 0.208092485549133
    "L11" = 0.5028901734104047
    "S11" = 0.008670520231213872
    "G12" = 0.7567567567567568
    "J12" = 0.1837837837837838
    "L12" = 0.005405405405405406
    "S12" = 0.05405405405405406
    "G13" = 0.6052631578947368
    "J13" = 0.2894736842105263
    "S13" = 0.1052631578947368
    "F15" = 0.02392344497607655
    "H15" = 0.1531100478468899
    "I15" = 0.05263157894736842
    "J15" = 0.3444976076555024
    "K15" = 0.01913875598086124
    "M15" = 0.009569377990430622
    "O15" = 0.06698564593301436
    "S15" = 0.3301435406698565
    "F16" = 0.05113636363636364
    "H16" = 0.1590909090909091
    "I16" = 0.05113636363636364
    "J16" = 0.3977272727272727
    "K16" = 0.1534090909090909
    "M16" = 0.005681818181818182
    "O16" = 0.04545454545454546
    "S16" = 0.1363636363636364
    "F17" = 0.03266331658291458
    "H17" = 0.1658291457286432
    "I17" = 0.1055276381909548
    "J17" = 0.407035175879397
    "K17" = 0.09547738693467336
    "M17" = 0.01758793969849246
    "N17" = 0.002512562814070352
    "O17" = 0.05276381909547739
    "S17" = 0.1206030150753769
    "F18" = 0.01092896174863388
    "H18" = 0.1311475409836066
    "I18" = 0.1147540983606557
    "J18" = 0.4262295081967213
    "K18" = 0.1147540983606557
    "M18" = 0.0273224043715847
    "N18" = 0.00546448087431694
    "O18" = 0.04918032786885246
    "S18" = 0.1202185792349727
    "F19" = 0.02579666160849772
    "H19" = 0.1965098634294385
    "I19" = 0.07814871016691957
    "J19" = 0.3566009104704097
    "K19" = 0.125948406676783
    "M19" = 0.01820940819423369
    "N19" = 0.0007587253414264037
    "O19" = 0.0629742033383915
    "S19" = 0.1350531107738998
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).Value = $changes[$cell]
}
